$wb = $excel.ActiveWorkbook

$wsWeb = $wb.Worksheets.Item("WEB_UI")
$wsWin = $wb.Worksheets.Item("WIN_UI")

# WEB_UI: clear D2 entirely (ClearContents removes the value)
$wsWeb.Range("D2").ClearContents()

# WIN_UI: set E2 to a single space (matches the blanked cells E3:E7)
$wsWin.Range("E2").Value = " "

# Restore selections as per the diff
[void]$wsWeb.Range("D2").Select()
[void]$wsWin.Range("E2").Select()
